$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Insert()

$ws.Range("A87").Value = 6
$ws.Range("B87").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44413
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = 100112032
$ws.Range("G87").Value = "Zapallo italiano"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 840
$ws.Range("K87").Value = 7000
$ws.Range("L87").Value = 8000
$ws.Range("M87").Value = 7536
$ws.Range("N87").Value = "`$/caja 50 unidades"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 151
$ws.Range("Q87").Value = 50
$ws.Range("R87").Value = "Hortaliza"
